$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final reconciled seminar/teacher-join data (nazev, zkratka, seminariciUcitIdno)
$data = @(
    ,@("Podnikové informační systémy", "EIS", 14)
    ,@("Podnikové informační systémy", "KEIS", 14)
    ,@("Fyzikální praktikum C", "K521", 302)
    ,@("Diplomový seminář", "K505", 306)
    ,@("Regional geography of the Czech Republic", "E101", 313)
    ,@("Reg. geography of Northwestern Bohemia", "E100", 313)
    ,@("Počítačové modelování I", "P107", 612)
    ,@("Programování A", "K103", 612)
    ,@("Počítačové modelování I", "K107", 612)
    ,@("Podnikové informační systémy", "EIS", 1609)
    ,@("Podnikové informační systémy", "KEIS", 1609)
    ,@("Identif. a hodn. ekosystémových služeb", "0153", 2527)
    ,@("Podnikové informační systémy", "EIS", 3457)
    ,@("Podnikové informační systémy", "KEIS", 3457)
    ,@("Podnikové informační systémy", "EIS", 3606)
    ,@("Podnikové informační systémy", "KEIS", 3606)
    ,@("Sociální sítě", "SON", 4190)
    ,@("Matematika I", "K106", 4221)
    ,@("Identif. a hodn. ekosystémových služeb", "0153", 4625)
    ,@("Praktické aplikace hardwaru", "AHW", 4746)
    ,@("Základy autonomní robotiky", "0182", 4746)
    ,@("Podnikové informační systémy", "EIS", 4991)
    ,@("Podnikové informační systémy", "KEIS", 4991)
    ,@("Reflektivní seminář pedagogické praxe", "KRSPP", 8021)
    ,@("Reflektivní seminář pedagogické praxe", "RSPP", 8021)
    ,@("Reflektivní seminář pedagogické praxe", "KSPP", 8021)
    ,@("Softwarové inženýrství", "SWI", 8093)
    ,@("Softwarové inženýrství", "KSWI", 8093)
    ,@("Odborná prezentace", "KOPRE", 8514)
    ,@("Odborná prezentace", "OPRE", 8514)
    ,@("Introduction to MATLAB", "ITM", 8514)
)

$rowCount = $data.Count
$startRow = 2

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    $ws.Cells.Item($r, 1).Value = $rowVals[0]
    $ws.Cells.Item($r, 2).Value = $rowVals[1]
    $ws.Cells.Item($r, 3).Value = $rowVals[2]
}

$lastRow = $startRow + $rowCount - 1

# Apply the same formatting used by the existing data rows (s="1" name/code cols, s="2" numeric col)
$ws.Range("A2:B" + $lastRow).VerticalAlignment = -4108
$numRange = $ws.Range("C2:C" + $lastRow)
$numRange.VerticalAlignment = -4108
$numRange.NumberFormat = "#,##0;[Red]-#,##0"

# Resize the table/autofilter to cover the new data extent
$lo = $ws.ListObjects.Item(1)
$newRange = $ws.Range("A1:C" + $lastRow)
$lo.Resize($newRange)
